$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"
$ws.Range("I14").Value = "sd"
$ws.Range("J14").Value = "Statement-non-opinion"
$ws.Range("I24").Value = "ba"
$ws.Range("J24").Value = "Appreciation"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I47").Value = "b"
$ws.Range("J47").Value = "Acknowledge (Backchannel)"
$ws.Range("I49").Value = "ba"
$ws.Range("J49").Value = "Appreciation"
$ws.Range("I55").Value = "sv"
$ws.Range("J55").Value = "Statement-opinion"
$ws.Range("I64").Value = "aa"
$ws.Range("J64").Value = "Agree/Accept"
$ws.Range("I65").Value = "%"
$ws.Range("J65").Value = "Uninterpretable"
$ws.Range("I68").Value = "b"
$ws.Range("J68").Value = "Acknowledge (Backchannel)"
$ws.Range("I79").Value = "aa"
$ws.Range("J79").Value = "Agree/Accept"
$ws.Range("I88").Value = "ba"
$ws.Range("J88").Value = "Appreciation"
$ws.Range("I97").Value = "b"
$ws.Range("J97").Value = "Acknowledge (Backchannel)"
$ws.Range("I99").Value = "sd"
$ws.Range("J99").Value = "Statement-non-opinion"
$ws.Range("I109").Value = "b"
$ws.Range("J109").Value = "Acknowledge (Backchannel)"
$ws.Range("I122").Value = "%"
$ws.Range("J122").Value = "Uninterpretable"
$ws.Range("I124").Value = "sd"
$ws.Range("J124").Value = "Statement-non-opinion"
$ws.Range("I131").Value = "sd"
$ws.Range("J131").Value = "Statement-non-opinion"
$ws.Range("I133").Value = "sv"
$ws.Range("J133").Value = "Statement-opinion"
$ws.Range("I135").Value = "sd"
$ws.Range("J135").Value = "Statement-non-opinion"
$ws.Range("I136").Value = "sd"
$ws.Range("J136").Value = "Statement-non-opinion"
$ws.Range("I138").Value = "aa"
$ws.Range("J138").Value = "Agree/Accept"
$ws.Range("I140").Value = "sd"
$ws.Range("J140").Value = "Statement-non-opinion"
$ws.Range("I153").Value = "aa"
$ws.Range("J153").Value = "Agree/Accept"
$ws.Range("I166").Value = "b"
$ws.Range("J166").Value = "Acknowledge (Backchannel)"
$ws.Range("I177").Value = "sd"
$ws.Range("J177").Value = "Statement-non-opinion"
$ws.Range("I184").Value = "ba"
$ws.Range("J184").Value = "Appreciation"
$ws.Range("I193").Value = "sd"
$ws.Range("J193").Value = "Statement-non-opinion"
$ws.Range("I201").Value = "b"
$ws.Range("J201").Value = "Acknowledge (Backchannel)"
$ws.Range("I202").Value = "ba"
$ws.Range("J202").Value = "Appreciation"
$ws.Range("I223").Value = "sd"
$ws.Range("J223").Value = "Statement-non-opinion"
$ws.Range("I232").Value = "%"
$ws.Range("J232").Value = "Uninterpretable"
$ws.Range("I259").Value = "sv"
$ws.Range("J259").Value = "Statement-opinion"
$ws.Range("I300").Value = "sd"
$ws.Range("J300").Value = "Statement-non-opinion"
$ws.Range("I302").Value = "aa"
$ws.Range("J302").Value = "Agree/Accept"
$ws.Range("I306").Value = "b"
$ws.Range("J306").Value = "Acknowledge (Backchannel)"
$ws.Range("I307").Value = "sv"
$ws.Range("J307").Value = "Statement-opinion"
$ws.Range("I320").Value = "sd"
$ws.Range("J320").Value = "Statement-non-opinion"
$ws.Range("I322").Value = "ba"
$ws.Range("J322").Value = "Appreciation"
$ws.Range("I323").Value = "aa"
$ws.Range("J323").Value = "Agree/Accept"
$ws.Range("I333").Value = "aa"
$ws.Range("J333").Value = "Agree/Accept"
$ws.Range("I341").Value = "sd"
$ws.Range("J341").Value = "Statement-non-opinion"
$ws.Range("I342").Value = "sv"
$ws.Range("J342").Value = "Statement-opinion"
$ws.Range("I345").Value = "aa"
$ws.Range("J345").Value = "Agree/Accept"
$ws.Range("I346").Value = "aa"
$ws.Range("J346").Value = "Agree/Accept"
$ws.Range("I357").Value = "sd"
$ws.Range("J357").Value = "Statement-non-opinion"
$ws.Range("I369").Value = "b"
$ws.Range("J369").Value = "Acknowledge (Backchannel)"
$ws.Range("I375").Value = "aa"
$ws.Range("J375").Value = "Agree/Accept"
$ws.Range("I376").Value = "sd"
$ws.Range("J376").Value = "Statement-non-opinion"
$ws.Range("I392").Value = "sd"
$ws.Range("J392").Value = "Statement-non-opinion"
$ws.Range("I403").Value = "b"
$ws.Range("J403").Value = "Acknowledge (Backchannel)"
$ws.Range("I405").Value = "sd"
$ws.Range("J405").Value = "Statement-non-opinion"
$ws.Range("I410").Value = "sd"
$ws.Range("J410").Value = "Statement-non-opinion"
$ws.Range("I413").Value = "sv"
$ws.Range("J413").Value = "Statement-opinion"
$ws.Range("I424").Value = "sd"
$ws.Range("J424").Value = "Statement-non-opinion"
$ws.Range("I444").Value = "aa"
$ws.Range("J444").Value = "Agree/Accept"
$ws.Range("I445").Value = "b"
$ws.Range("J445").Value = "Acknowledge (Backchannel)"
$ws.Range("I447").Value = "%"
$ws.Range("J447").Value = "Uninterpretable"
$ws.Range("I452").Value = "ba"
$ws.Range("J452").Value = "Appreciation"
$ws.Range("I457").Value = "ba"
$ws.Range("J457").Value = "Appreciation"
$ws.Range("I476").Value = "sd"
$ws.Range("J476").Value = "Statement-non-opinion"
$ws.Range("I477").Value = "sd"
$ws.Range("J477").Value = "Statement-non-opinion"
$ws.Range("I484").Value = "sd"
$ws.Range("J484").Value = "Statement-non-opinion"
$ws.Range("I492").Value = "sd"
$ws.Range("J492").Value = "Statement-non-opinion"
